{"js": "const replacements = [\n  [\"2025-04-23 Wednesday\", \"2025-04-24 Thursday\"],\n  [\"791\u00f73=\", \"773\u00f77=\"],\n  [\"742\u00f76=\", \"127\u00f76=\"],\n  [\"367\u00f75=\", \"577\u00f79=\"],\n  [\"536\u00f76=\", \"366\u00f78=\"],\n  [\"109\u00f78=\", \"457\u00f73=\"],\n  [\"211\u00f76=\", \"554\u00f72=\"],\n  [\"376\u00f78=\", \"601\u00f76=\"],\n  [\"907\u00f76=\", \"492\u00f77=\"],\n  [\"850\u00f72=\", \"533\u00f75=\"],\n  [\"447\u00f74=\", \"688\u00f74=\"],\n  [\"297\u00f74=\", \"675\u00f75=\"],\n  [\"360\u00f76=\", \"779\u00f76=\"],\n  [\"554\u00f77=\", \"458\u00f72=\"],\n  [\"387\u00f77=\", \"336\u00f75=\"],\n  [\"342\u00f77=\", \"268\u00f76=\"],\n  [\"161\u00f76=\", \"858\u00f72=\"],\n  [\"565\u00f75=\", \"878\u00f77=\"],\n  [\"220\u00f73=\", \"701\u00f73=\"],\n  [\"139\u00f79=\", \"751\u00f72=\"],\n  [\"434\u00f74=\", \"849\u00f75=\"],\n  [\"895\u00f76=\", \"227\u00f78=\"],\n  [\"508\u00f79=\", \"162\u00f78=\"],\n  [\"165\u00f79=\", \"630\u00f72=\"],\n  [\"365\u00f78=\", \"804\u00f73=\"],\n  [\"270\u00f72=\", \"713\u00f79=\"],\n];\n\nfor (const [searchText, replaceText] of replacements) {\n  const results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"2025-04-23 Wednesday\"; Replace = \"2025-04-24 Thursday\" },\n    @{ Find = \"791\u00f73=\"; Replace = \"773\u00f77=\" },\n    @{ Find = \"742\u00f76=\"; Replace = \"127\u00f76=\" },\n    @{ Find = \"367\u00f75=\"; Replace = \"577\u00f79=\" },\n    @{ Find = \"536\u00f76=\"; Replace = \"366\u00f78=\" },\n    @{ Find = \"109\u00f78=\"; Replace = \"457\u00f73=\" },\n    @{ Find = \"211\u00f76=\"; Replace = \"554\u00f72=\" },\n    @{ Find = \"376\u00f78=\"; Replace = \"601\u00f76=\" },\n    @{ Find = \"907\u00f76=\"; Replace = \"492\u00f77=\" },\n    @{ Find = \"850\u00f72=\"; Replace = \"533\u00f75=\" },\n    @{ Find = \"447\u00f74=\"; Replace = \"688\u00f74=\" },\n    @{ Find = \"297\u00f74=\"; Replace = \"675\u00f75=\" },\n    @{ Find = \"360\u00f76=\"; Replace = \"779\u00f76=\" },\n    @{ Find = \"554\u00f77=\"; Replace = \"458\u00f72=\" },\n    @{ Find = \"387\u00f77=\"; Replace = \"336\u00f75=\" },\n    @{ Find = \"342\u00f77=\"; Replace = \"268\u00f76=\" },\n    @{ Find = \"161\u00f76=\"; Replace = \"858\u00f72=\" },\n    @{ Find = \"565\u00f75=\"; Replace = \"878\u00f77=\" },\n    @{ Find = \"220\u00f73=\"; Replace = \"701\u00f73=\" },\n    @{ Find = \"139\u00f79=\"; Replace = \"751\u00f72=\" },\n    @{ Find = \"434\u00f74=\"; Replace = \"849\u00f75=\" },\n    @{ Find = \"895\u00f76=\"; Replace = \"227\u00f78=\" },\n    @{ Find = \"508\u00f79=\"; Replace = \"162\u00f78=\" },\n    @{ Find = \"165\u00f79=\"; Replace = \"630\u00f72=\" },\n    @{ Find = \"365\u00f78=\"; Replace = \"804\u00f73=\" },\n    @{ Find = \"270\u00f72=\"; Replace = \"713\u00f79=\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.Replace\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute(\n        $r.Find,      # FindText\n        $true,        # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        1,            # Wrap (wdFindContinue)\n        $false,       # Format\n        $r.Replace,   # ReplaceWith\n        2             # Replace (wdReplaceAll)\n    )\n}\n"}
